# Applies the "fixed order, and phi distribution of particles" edit:
#  - Number of Events (D4): 2000 -> 5000
#  - New per-system timing values in column D (rows 7-21)
#  - Row 15/16 labels swapped (torus/dc order fixed)
#  - Formulas in J column simplified to divide by $D$21 instead of $I$21
#  - Selection moved to J12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Number of events -------------------------------------------------
$ws.Range("D4").Value = 5000

# ---- New label/value pairs for rows 7-21 -------------------------------
# Each entry: row, label (column C), new raw time value (column D)
$rowsData = @(
    @(7,  "target",   5.9450000000000003),
    @(8,  "svt",      10.420299999999999),
    @(9,  "ctof",     11.6716),
    @(10, "cnd",      12.524800000000001),
    @(11, "solenoid", 181.00299999999999),
    @(12, "mm",       288.33),
    @(13, "htcc",     414.76299999999998),
    @(14, "ft",       624.875),
    @(15, "dc",       579.625),
    @(16, "torus",    594.04499999999996),
    @(17, "ltcc",     590.26400000000001),
    @(18, "rich",     603.63),
    @(19, "ftof",     624.875),
    @(20, "pcal",     1586.6),
    @(21, "ecAll",    2587.17)
)

foreach ($item in $rowsData) {
    $r = $item[0]
    $label = $item[1]
    $val = $item[2]

    $ws.Range("C$r").Value = $label
    $ws.Range("D$r").Value = $val
}

# ---- Formulas -----------------------------------------------------------
# E, F, G, H, I keep the same per-row formula shapes; only row 7/8 anchor
# the "D7"/"D8-D7" seed, the rest follow the same relative pattern.
for ($r = 7; $r -le 21; $r++) {
    $ws.Range("E$r").Formula = '=1000*D' + $r + '/$D$4'
    $ws.Range("F$r").Formula = '=$D$4/D' + $r
    $ws.Range("G$r").Formula = '=D' + $r + '/D$21'
    $ws.Range("H$r").Formula = '=C' + $r
}

$ws.Range("I7").Formula = '=D7'
for ($r = 8; $r -le 21; $r++) {
    $prev = $r - 1
    $ws.Range("I$r").Formula = '=D' + $r + '-D' + $prev
}

# J now divides by the grand total in $D$21 rather than $I$21
for ($r = 7; $r -le 21; $r++) {
    $ws.Range("J$r").Formula = '=I' + $r + '/$D$21'
}

# ---- Selection ------------------------------------------------------------
$ws.Range("J12").Select()
